$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 84.07754166666668
$ws.Range("H2").Value = 252.232625
$ws.Range("I2").Value = 0.5195714800795683
$ws.Range("J2").Value = 0.5195714800795683
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 84.07754166666668
$ws.Range("N2").Value = 252.232625
$ws.Range("O2").Value = 0.5195714800795683
$ws.Range("P2").Value = 0.5195714800795683
$ws.Range("Q2").Value = 7069.033012710071
$ws.Range("R2").Value = 63621.29711439063
$ws.Range("S2").Value = 0.2699545229120733
$ws.Range("T2").Value = 0.2699545229120733

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 84.07754166666668
$ws.Range("H3").Value = 252.232625
$ws.Range("I3").Value = 0.5195714800795683
$ws.Range("J3").Value = 0.5195714800795683
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.717126
$ws.Range("N3").Value = 2.151378
$ws.Range("O3").Value = 0.004431602183383778
$ws.Range("P3").Value = 0.004431602183383778
$ws.Range("Q3").Value = 60.29419114525001
$ws.Range("R3").Value = 542.6477203072501
$ws.Range("S3").Value = 0.002302534105544556
$ws.Range("T3").Value = 0.002302534105544556

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 84.07754166666668
$ws.Range("H4").Value = 252.232625
$ws.Range("I4").Value = 0.5195714800795683
$ws.Range("J4").Value = 0.5195714800795683
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 77.02626533333334
$ws.Range("N4").Value = 231.078796
$ws.Range("O4").Value = 0.4759969177370478
$ws.Range("P4").Value = 0.4759969177370479
$ws.Range("Q4").Value = 6476.179032991057
$ws.Range("R4").Value = 58285.61129691951
$ws.Range("S4").Value = 0.2473144230619505
$ws.Range("T4").Value = 0.2473144230619505

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.717126
$ws.Range("H5").Value = 2.151378
$ws.Range("I5").Value = 0.004431602183383778
$ws.Range("J5").Value = 0.004431602183383778
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 84.07754166666668
$ws.Range("N5").Value = 252.232625
$ws.Range("O5").Value = 0.5195714800795683
$ws.Range("P5").Value = 0.5195714800795683
$ws.Range("Q5").Value = 60.29419114525001
$ws.Range("R5").Value = 542.6477203072501
$ws.Range("S5").Value = 0.002302534105544556
$ws.Range("T5").Value = 0.002302534105544556

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.717126
$ws.Range("H6").Value = 2.151378
$ws.Range("I6").Value = 0.004431602183383778
$ws.Range("J6").Value = 0.004431602183383778
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.717126
$ws.Range("N6").Value = 2.151378
$ws.Range("O6").Value = 0.004431602183383778
$ws.Range("P6").Value = 0.004431602183383778
$ws.Range("Q6").Value = 0.514269699876
$ws.Range("R6").Value = 4.628427298884001
$ws.Range("S6").Value = 0.00001963909791177186
$ws.Range("T6").Value = 0.00001963909791177187

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.717126
$ws.Range("H7").Value = 2.151378
$ws.Range("I7").Value = 0.004431602183383778
$ws.Range("J7").Value = 0.004431602183383778
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 77.02626533333334
$ws.Range("N7").Value = 231.078796
$ws.Range("O7").Value = 0.4759969177370478
$ws.Range("P7").Value = 0.4759969177370479
$ws.Range("Q7").Value = 55.23753755343201
$ws.Range("R7").Value = 497.1378379808881
$ws.Range("S7").Value = 0.00210942897992745
$ws.Range("T7").Value = 0.00210942897992745

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 77.02626533333334
$ws.Range("H8").Value = 231.078796
$ws.Range("I8").Value = 0.4759969177370478
$ws.Range("J8").Value = 0.4759969177370479
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 84.07754166666668
$ws.Range("N8").Value = 252.232625
$ws.Range("O8").Value = 0.5195714800795683
$ws.Range("P8").Value = 0.5195714800795683
$ws.Range("Q8").Value = 6476.179032991057
$ws.Range("R8").Value = 58285.61129691951
$ws.Range("S8").Value = 0.2473144230619505
$ws.Range("T8").Value = 0.2473144230619505

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 77.02626533333334
$ws.Range("H9").Value = 231.078796
$ws.Range("I9").Value = 0.4759969177370478
$ws.Range("J9").Value = 0.4759969177370479
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.717126
$ws.Range("N9").Value = 2.151378
$ws.Range("O9").Value = 0.004431602183383778
$ws.Range("P9").Value = 0.004431602183383778
$ws.Range("Q9").Value = 55.23753755343201
$ws.Range("R9").Value = 497.1378379808881
$ws.Range("S9").Value = 0.00210942897992745
$ws.Range("T9").Value = 0.00210942897992745

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 77.02626533333334
$ws.Range("H10").Value = 231.078796
$ws.Range("I10").Value = 0.4759969177370478
$ws.Range("J10").Value = 0.4759969177370479
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 77.02626533333334
$ws.Range("N10").Value = 231.078796
$ws.Range("O10").Value = 0.4759969177370478
$ws.Range("P10").Value = 0.4759969177370479
$ws.Range("Q10").Value = 5933.045551201069
$ws.Range("R10").Value = 53397.40996080962
$ws.Range("S10").Value = 0.2265730656951699
$ws.Range("T10").Value = 0.2265730656951699

